# Auto-generated script applying numeric updates to the Siren_Profits workbook
# (commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 640.0769
$ws.Range("I2").Value = 651.75
$ws.Range("K2").Value = 651.75
$ws.Range("M2").Value = -538.75
$ws.Range("H11").Value = 135.5
$ws.Range("I11").Value = 135.5
$ws.Range("K11").Value = 135.5
$ws.Range("M11").Value = 4.5
$ws.Range("H15").Value = 2209.3438
$ws.Range("I15").Value = 2209.3438
$ws.Range("K15").Value = 6628.0314
$ws.Range("M15").Value = -6459.0314
$ws.Range("H18").Value = 6833.7
$ws.Range("I18").Value = 8399.286
$ws.Range("K18").Value = 8399.286
$ws.Range("M18").Value = -8115.286
$ws.Range("H19").Value = 1898.5
$ws.Range("I19").Value = 997
$ws.Range("K19").Value = 997
$ws.Range("M19").Value = -822
$ws.Range("H26").Value = 19998.5
$ws.Range("J26").Value = 19998.5
$ws.Range("L26").Value = 19998.5
$ws.Range("N26").Value = -20686.5
$ws.Range("H32").Value = 5951.857
$ws.Range("J32").Value = 5897.6665
$ws.Range("L32").Value = 5897.6665
$ws.Range("N32").Value = -6549.6665
$ws.Range("H40").Value = 4998
$ws.Range("H42").Value = 697.5714
$ws.Range("I42").Value = 93
$ws.Range("J42").Value = 939.4
$ws.Range("K42").Value = 279
$ws.Range("L42").Value = 2818.2
$ws.Range("M42").Value = -49
$ws.Range("N42").Value = -3278.2
$ws.Range("H43").Value = 6915.3335
$ws.Range("J43").Value = 7623
$ws.Range("L43").Value = 7623
$ws.Range("N43").Value = -7761
$ws.Range("H46").Value = 4750
$ws.Range("J46").Value = 4750
$ws.Range("L46").Value = 14250
$ws.Range("N46").Value = -14488
$ws.Range("H51").Value = 7379
$ws.Range("I51").Value = 6857.8
$ws.Range("J51").Value = 9985
$ws.Range("K51").Value = 6857.8
$ws.Range("L51").Value = 9985
$ws.Range("M51").Value = -6373.8
$ws.Range("N51").Value = -10953
$ws.Range("H53").Value = 5308.5
$ws.Range("I53").Value = 6981.6665
$ws.Range("J53").Value = 289
$ws.Range("K53").Value = 6981.6665
$ws.Range("L53").Value = 289
$ws.Range("M53").Value = -6344.6665
$ws.Range("N53").Value = -1563
$ws.Range("H60").Value = 4750
$ws.Range("J60").Value = 4750
$ws.Range("L60").Value = 14250
$ws.Range("N60").Value = -15218
$ws.Range("H62").Value = 20205.105
$ws.Range("I62").Value = 6612.125
$ws.Range("K62").Value = 6612.125
$ws.Range("M62").Value = -5988.125
$ws.Range("H65").Value = 20205.105
$ws.Range("I65").Value = 6612.125
$ws.Range("K65").Value = 33060.625
$ws.Range("M65").Value = -29940.625
$ws.Range("H75").Value = 123000
$ws.Range("J75").Value = 123000
$ws.Range("L75").Value = 123000
$ws.Range("N75").Value = -124872
$ws.Range("H78").Value = 123000
$ws.Range("J78").Value = 123000
$ws.Range("L78").Value = 369000
$ws.Range("N78").Value = -378360
$ws.Range("H86").Value = 47264.8
$ws.Range("J86").Value = 91799.8
$ws.Range("L86").Value = 91799.8
$ws.Range("N86").Value = -94045.8
$ws.Range("H89").Value = 47264.8
$ws.Range("J89").Value = 91799.8
$ws.Range("L89").Value = 458999
$ws.Range("N89").Value = -470231
$ws.Range("H92").Value = 1396.8
$ws.Range("I92").Value = 969.3333
$ws.Range("J92").Value = 1580
$ws.Range("K92").Value = 969.3333
$ws.Range("L92").Value = 1580
$ws.Range("M92").Value = 278.6667
$ws.Range("N92").Value = -4076
$ws.Range("H107").Value = 9080.895
$ws.Range("I107").Value = 10357.154
$ws.Range("J107").Value = 6315.6665
$ws.Range("K107").Value = 10357.154
$ws.Range("L107").Value = 6315.6665
$ws.Range("M107").Value = -8437.154
$ws.Range("N107").Value = -10155.6665
$ws.Range("H109").Value = 456025000
$ws.Range("J109").Value = 456025000
$ws.Range("L109").Value = 456025000
$ws.Range("N109").Value = -456027774
$ws.Range("H112").Value = 1459.8857
$ws.Range("J112").Value = 1517.1538
$ws.Range("L112").Value = 4551.4614
$ws.Range("N112").Value = -6767.4614
$ws.Range("H113").Value = 8309.700000000001
$ws.Range("I113").Value = 8168.375
$ws.Range("J113").Value = 8875
$ws.Range("K113").Value = 8168.375
$ws.Range("L113").Value = 8875
$ws.Range("M113").Value = -4914.375
$ws.Range("N113").Value = -15383
$ws.Range("H135").Value = 21999.75
$ws.Range("J135").Value = 4000
$ws.Range("L135").Value = 36000
$ws.Range("N135").Value = -41070
$ws.Range("H137").Value = 7740.75
$ws.Range("I137").Value = 9066.893
$ws.Range("J137").Value = 3099.25
$ws.Range("K137").Value = 27200.679
$ws.Range("L137").Value = 9297.75
$ws.Range("M137").Value = -24650.679
$ws.Range("N137").Value = -14397.75
$ws.Range("H138").Value = 3345.814
$ws.Range("I138").Value = 720.6429000000001
$ws.Range("J138").Value = 4613.1377
$ws.Range("K138").Value = 2161.9287
$ws.Range("L138").Value = 13839.4131
$ws.Range("M138").Value = 2978.0713
$ws.Range("N138").Value = -24119.4131
$ws.Range("H140").Value = 80697.60000000001
$ws.Range("J140").Value = 80697.60000000001
$ws.Range("L140").Value = 80697.60000000001
$ws.Range("N140").Value = -91057.60000000001
$ws.Range("H141").Value = 4999.357
$ws.Range("I141").Value = 4208.9
$ws.Range("J141").Value = 6975.5
$ws.Range("K141").Value = 12626.7
$ws.Range("L141").Value = 20926.5
$ws.Range("M141").Value = -7446.699999999999
$ws.Range("N141").Value = -31286.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 4950
$ws.Range("H32").Value = 1766.86
$ws.Range("I32").Value = 1762.102
$ws.Range("K32").Value = 1762.102
$ws.Range("M32").Value = -1475.102
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H45").Value = 9737.375
$ws.Range("I45").Value = 11800
$ws.Range("K45").Value = 11800
$ws.Range("M45").Value = -11423
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 4457.9473
$ws.Range("I61").Value = 3694.4443
$ws.Range("J61").Value = 6332
$ws.Range("K61").Value = 3694.4443
$ws.Range("L61").Value = 6332
$ws.Range("M61").Value = -3482.4443
$ws.Range("N61").Value = -6756
$ws.Range("H63").Value = 4081.1667
$ws.Range("I63").Value = 4498.2
$ws.Range("J63").Value = 1996
$ws.Range("K63").Value = 4498.2
$ws.Range("L63").Value = 1996
$ws.Range("M63").Value = -3812.2
$ws.Range("N63").Value = -3368
$ws.Range("H66").Value = 4081.1667
$ws.Range("I66").Value = 4498.2
$ws.Range("J66").Value = 1996
$ws.Range("K66").Value = 22491
$ws.Range("L66").Value = 9980
$ws.Range("M66").Value = -19059
$ws.Range("N66").Value = -16844
$ws.Range("H102").Value = 27677.111
$ws.Range("J102").Value = 8819
$ws.Range("L102").Value = 8819
$ws.Range("N102").Value = -12063
$ws.Range("H110").Value = 8486.058999999999
$ws.Range("I110").Value = 11296.728
$ws.Range("J110").Value = 3333.1667
$ws.Range("K110").Value = 11296.728
$ws.Range("L110").Value = 3333.1667
$ws.Range("M110").Value = -9251.727999999999
$ws.Range("N110").Value = -7423.1667
$ws.Range("H122").Value = 615351.25
$ws.Range("I122").Value = 6286
$ws.Range("J122").Value = 940186.0600000001
$ws.Range("K122").Value = 18858
$ws.Range("L122").Value = 2820558.18
$ws.Range("M122").Value = -16408
$ws.Range("N122").Value = -2825458.18
$ws.Range("H132").Value = 4621.7666
$ws.Range("I132").Value = 3738.2632
$ws.Range("J132").Value = 6147.8184
$ws.Range("K132").Value = 11214.7896
$ws.Range("L132").Value = 18443.4552
$ws.Range("M132").Value = -8684.7896
$ws.Range("N132").Value = -23503.4552
$ws.Range("H136").Value = 4457.9473
$ws.Range("I136").Value = 3694.4443
$ws.Range("J136").Value = 6332
$ws.Range("K136").Value = 11083.3329
$ws.Range("L136").Value = 18996
$ws.Range("M136").Value = -8533.332900000001
$ws.Range("N136").Value = -24096
$ws.Range("H138").Value = 74421.78
$ws.Range("J138").Value = 74421.78
$ws.Range("L138").Value = 74421.78
$ws.Range("N138").Value = -84701.78

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5566
$ws.Range("H81").Value = 28666.666
$ws.Range("J81").Value = 28666.666
$ws.Range("L81").Value = 28666.666
$ws.Range("N81").Value = -30788.666
$ws.Range("H84").Value = 28666.666
$ws.Range("J84").Value = 28666.666
$ws.Range("L84").Value = 85999.99800000001
$ws.Range("N84").Value = -96607.99800000001
$ws.Range("H86").Value = 8382.333000000001
$ws.Range("I86").Value = 12616.833
$ws.Range("K86").Value = 12616.833
$ws.Range("M86").Value = -11493.833
$ws.Range("H89").Value = 8382.333000000001
$ws.Range("I89").Value = 12616.833
$ws.Range("K89").Value = 63084.165
$ws.Range("M89").Value = -57468.165
$ws.Range("H105").Value = 4966.154
$ws.Range("I105").Value = 2695
$ws.Range("J105").Value = 8600
$ws.Range("K105").Value = 2695
$ws.Range("L105").Value = 8600
$ws.Range("M105").Value = -948
$ws.Range("N105").Value = -12094
$ws.Range("H107").Value = 11601.929
$ws.Range("I107").Value = 11725.154
$ws.Range("K107").Value = 11725.154
$ws.Range("M107").Value = -9805.154

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H22").Value = 1358.5
$ws.Range("I22").Value = 825
$ws.Range("J22").Value = 1587.1428
$ws.Range("K22").Value = 825
$ws.Range("L22").Value = 1587.1428
$ws.Range("M22").Value = -475
$ws.Range("N22").Value = -2287.1428
$ws.Range("H25").Value = 1013
$ws.Range("J25").Value = 1013
$ws.Range("L25").Value = 1013
$ws.Range("N25").Value = -1361
$ws.Range("H31").Value = 3515
$ws.Range("I31").Value = 1200.8572
$ws.Range("K31").Value = 1200.8572
$ws.Range("M31").Value = -905.8571999999999
$ws.Range("H34").Value = 3515
$ws.Range("I34").Value = 1200.8572
$ws.Range("K34").Value = 1200.8572
$ws.Range("M34").Value = -998.8571999999999
$ws.Range("H58").Value = 1859.875
$ws.Range("I58").Value = 673.6923
$ws.Range("K58").Value = 673.6923
$ws.Range("M58").Value = -470.6923
$ws.Range("H62").Value = 60892.547
$ws.Range("I62").Value = 8619.666999999999
$ws.Range("J62").Value = 123620
$ws.Range("K62").Value = 8619.666999999999
$ws.Range("L62").Value = 123620
$ws.Range("M62").Value = -7995.666999999999
$ws.Range("N62").Value = -124868
$ws.Range("H65").Value = 60892.547
$ws.Range("I65").Value = 8619.666999999999
$ws.Range("J65").Value = 123620
$ws.Range("K65").Value = 43098.335
$ws.Range("L65").Value = 618100
$ws.Range("M65").Value = -39978.335
$ws.Range("N65").Value = -624340
$ws.Range("H86").Value = 12911.29
$ws.Range("J86").Value = 13497.1
$ws.Range("L86").Value = 13497.1
$ws.Range("N86").Value = -15743.1
$ws.Range("H88").Value = 43329.332
$ws.Range("I88").Value = 50000
$ws.Range("J88").Value = 39994
$ws.Range("K88").Value = 50000
$ws.Range("L88").Value = 39994
$ws.Range("M88").Value = -49594
$ws.Range("N88").Value = -40806
$ws.Range("H89").Value = 12911.29
$ws.Range("J89").Value = 13497.1
$ws.Range("L89").Value = 67485.5
$ws.Range("N89").Value = -78717.5
$ws.Range("H91").Value = 43329.332
$ws.Range("I91").Value = 50000
$ws.Range("J91").Value = 39994
$ws.Range("K91").Value = 50000
$ws.Range("L91").Value = 39994
$ws.Range("M91").Value = -48596
$ws.Range("N91").Value = -42802
$ws.Range("H94").Value = 2236.5
$ws.Range("J94").Value = 1961.25
$ws.Range("L94").Value = 1961.25
$ws.Range("N94").Value = -2863.25
$ws.Range("H119").Value = 49999
$ws.Range("J119").Value = 49999
$ws.Range("L119").Value = 49999
$ws.Range("N119").Value = -59675
$ws.Range("H132").Value = 60375
$ws.Range("I132").Value = 5750
$ws.Range("K132").Value = 17250
$ws.Range("M132").Value = -14720
$ws.Range("H134").Value = 7639.143
$ws.Range("I134").Value = 7742
$ws.Range("J134").Value = 7502
$ws.Range("K134").Value = 23226
$ws.Range("L134").Value = 22506
$ws.Range("M134").Value = -20691
$ws.Range("N134").Value = -27576
$ws.Range("H136").Value = 1859.875
$ws.Range("I136").Value = 673.6923
$ws.Range("K136").Value = 2021.0769
$ws.Range("M136").Value = 528.9231
$ws.Range("H138").Value = 22999.5
$ws.Range("J138").Value = 22999.5
$ws.Range("L138").Value = 22999.5
$ws.Range("N138").Value = -33279.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 743.5
$ws.Range("I3").Value = 492.2
$ws.Range("K3").Value = 1476.6
$ws.Range("M3").Value = -1364.6
$ws.Range("H6").Value = 2122
$ws.Range("I6").Value = 1490.1428
$ws.Range("K6").Value = 4470.428400000001
$ws.Range("M6").Value = -4357.428400000001
$ws.Range("H33").Value = 121.57143
$ws.Range("I33").Value = 104.72727
$ws.Range("J33").Value = 183.33333
$ws.Range("K33").Value = 628.3636200000001
$ws.Range("L33").Value = 1099.99998
$ws.Range("M33").Value = -345.3636200000001
$ws.Range("N33").Value = -1665.99998
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H97").Value = 27918.564
$ws.Range("J97").Value = 1954.875
$ws.Range("L97").Value = 5864.625
$ws.Range("N97").Value = -6856.625
$ws.Range("H98").Value = 959.25
$ws.Range("I98").Value = 1081.1666
$ws.Range("J98").Value = 837.3333
$ws.Range("K98").Value = 3243.4998
$ws.Range("L98").Value = 2511.9999
$ws.Range("M98").Value = -1745.4998
$ws.Range("N98").Value = -5507.9999
$ws.Range("H99").Value = 6517.5454
$ws.Range("I99").Value = 1431.25
$ws.Range("J99").Value = 9424
$ws.Range("K99").Value = 4293.75
$ws.Range("L99").Value = 28272
$ws.Range("M99").Value = -2047.75
$ws.Range("N99").Value = -32764
$ws.Range("H115").Value = 9414.333000000001
$ws.Range("I115").Value = 9414.333000000001
$ws.Range("K115").Value = 28242.999
$ws.Range("M115").Value = -27067.999
$ws.Range("H122").Value = 2315.4138
$ws.Range("J122").Value = 3395.125
$ws.Range("L122").Value = 30556.125
$ws.Range("N122").Value = -35456.125
$ws.Range("H131").Value = 8793.6
$ws.Range("J131").Value = 2145.75
$ws.Range("L131").Value = 6437.25
$ws.Range("N131").Value = -16517.25
$ws.Range("H132").Value = 101377.2
$ws.Range("I132").Value = 962
$ws.Range("J132").Value = 252000
$ws.Range("K132").Value = 8658
$ws.Range("L132").Value = 2268000
$ws.Range("M132").Value = -6128
$ws.Range("N132").Value = -2273060
$ws.Range("H139").Value = 1501869.5
$ws.Range("I139").Value = 1668188.2
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 5004564.6
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -4999424.6
$ws.Range("N139").Value = -25280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 28328.166
$ws.Range("J57").Value = 32492.25
$ws.Range("L57").Value = 32492.25
$ws.Range("N57").Value = -34132.25
$ws.Range("H70").Value = 6447.3335
$ws.Range("I70").Value = 5761.5835
$ws.Range("J70").Value = 7361.6665
$ws.Range("K70").Value = 5761.5835
$ws.Range("L70").Value = 7361.6665
$ws.Range("M70").Value = -5491.5835
$ws.Range("N70").Value = -7901.6665
$ws.Range("H73").Value = 6447.3335
$ws.Range("I73").Value = 5761.5835
$ws.Range("J73").Value = 7361.6665
$ws.Range("K73").Value = 5761.5835
$ws.Range("L73").Value = 7361.6665
$ws.Range("M73").Value = -4825.5835
$ws.Range("N73").Value = -9233.666499999999
$ws.Range("H102").Value = 7765
$ws.Range("I102").Value = 10322
$ws.Range("K102").Value = 10322
$ws.Range("M102").Value = -8700
$ws.Range("H109").Value = 153134.25
$ws.Range("J109").Value = 153134.25
$ws.Range("L109").Value = 153134.25
$ws.Range("N109").Value = -155214.25
$ws.Range("H122").Value = 47440
$ws.Range("I122").Value = 50666.668
$ws.Range("J122").Value = 42600
$ws.Range("K122").Value = 152000.004
$ws.Range("L122").Value = 127800
$ws.Range("M122").Value = -149550.004
$ws.Range("N122").Value = -132700
$ws.Range("H132").Value = 4004.5715
$ws.Range("I132").Value = 3928
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11784
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -9254
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 95893.5
$ws.Range("J136").Value = 95893.5
$ws.Range("L136").Value = 287680.5
$ws.Range("N136").Value = -292780.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 29893.059
$ws.Range("I7").Value = 55585.75
$ws.Range("J7").Value = 7055.1113
$ws.Range("K7").Value = 55585.75
$ws.Range("L7").Value = 7055.1113
$ws.Range("M7").Value = -55473.75
$ws.Range("N7").Value = -7279.1113
$ws.Range("H16").Value = 4679.278
$ws.Range("J16").Value = 5950.6
$ws.Range("L16").Value = 5950.6
$ws.Range("N16").Value = -6290.6
$ws.Range("H22").Value = 3752.75
$ws.Range("I22").Value = 3781.4443
$ws.Range("J22").Value = 3666.6667
$ws.Range("K22").Value = 3781.4443
$ws.Range("L22").Value = 3666.6667
$ws.Range("M22").Value = -3486.4443
$ws.Range("N22").Value = -4256.6667
$ws.Range("H24").Value = 14083.667
$ws.Range("I24").Value = 14083.667
$ws.Range("K24").Value = 14083.667
$ws.Range("M24").Value = -13740.667
$ws.Range("H27").Value = 3752.75
$ws.Range("I27").Value = 3781.4443
$ws.Range("J27").Value = 3666.6667
$ws.Range("K27").Value = 3781.4443
$ws.Range("L27").Value = 3666.6667
$ws.Range("M27").Value = -3674.4443
$ws.Range("N27").Value = -3880.6667
$ws.Range("H46").Value = 3527.3125
$ws.Range("J46").Value = 5098.8
$ws.Range("L46").Value = 5098.8
$ws.Range("N46").Value = -5474.8
$ws.Range("H68").Value = 4455.6665
$ws.Range("I68").Value = 1975.5
$ws.Range("K68").Value = 1975.5
$ws.Range("M68").Value = -1226.5
$ws.Range("H71").Value = 4455.6665
$ws.Range("I71").Value = 1975.5
$ws.Range("K71").Value = 9877.5
$ws.Range("M71").Value = -6133.5
$ws.Range("H122").Value = 4028.68
$ws.Range("I122").Value = 3573.9375
$ws.Range("J122").Value = 4837.1113
$ws.Range("K122").Value = 10721.8125
$ws.Range("L122").Value = 14511.3339
$ws.Range("M122").Value = -8271.8125
$ws.Range("N122").Value = -19411.3339
$ws.Range("H126").Value = 29893.059
$ws.Range("I126").Value = 55585.75
$ws.Range("J126").Value = 7055.1113
$ws.Range("K126").Value = 166757.25
$ws.Range("L126").Value = 21165.3339
$ws.Range("M126").Value = -164287.25
$ws.Range("N126").Value = -26105.3339
$ws.Range("H132").Value = 375695.4
$ws.Range("I132").Value = 574985.4399999999
$ws.Range("K132").Value = 1724956.32
$ws.Range("M132").Value = -1722426.32
$ws.Range("H134").Value = 112946.664
$ws.Range("J134").Value = 112946.664
$ws.Range("L134").Value = 112946.664
$ws.Range("N134").Value = -123086.664
$ws.Range("H136").Value = 10006.523
$ws.Range("I136").Value = 27983.334
$ws.Range("K136").Value = 83950.00199999999
$ws.Range("M136").Value = -81400.00199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 92666.664
$ws.Range("J64").Value = 92666.664
$ws.Range("L64").Value = 92666.664
$ws.Range("N64").Value = -93162.664
$ws.Range("H67").Value = 92666.664
$ws.Range("J67").Value = 92666.664
$ws.Range("L67").Value = 92666.664
$ws.Range("N67").Value = -94382.664
$ws.Range("H81").Value = 11504.733
$ws.Range("I81").Value = 21709.285
$ws.Range("K81").Value = 43418.57
$ws.Range("M81").Value = -42357.57
$ws.Range("H82").Value = 44999.5
$ws.Range("J82").Value = 44999
$ws.Range("L82").Value = 44999
$ws.Range("N82").Value = -45765
$ws.Range("H84").Value = 11504.733
$ws.Range("I84").Value = 21709.285
$ws.Range("K84").Value = 217092.85
$ws.Range("M84").Value = -211788.85
$ws.Range("H85").Value = 44999.5
$ws.Range("J85").Value = 44999
$ws.Range("L85").Value = 44999
$ws.Range("N85").Value = -47651
$ws.Range("H122").Value = 21661.818
$ws.Range("I122").Value = 3920.9524
$ws.Range("J122").Value = 52708.332
$ws.Range("K122").Value = 11762.8572
$ws.Range("L122").Value = 158124.996
$ws.Range("M122").Value = -9312.8572
$ws.Range("N122").Value = -163024.996
$ws.Range("H126").Value = 29424.188
$ws.Range("I126").Value = 39027
$ws.Range("K126").Value = 117081
$ws.Range("M126").Value = -114611
$ws.Range("H128").Value = 45569.9
$ws.Range("J128").Value = 65249.5
$ws.Range("L128").Value = 65249.5
$ws.Range("N128").Value = -75209.5
$ws.Range("H136").Value = 2372.3
$ws.Range("I136").Value = 1713.1333
$ws.Range("J136").Value = 4349.8
$ws.Range("K136").Value = 5139.3999
$ws.Range("L136").Value = 13049.4
$ws.Range("M136").Value = -2589.3999
$ws.Range("N136").Value = -18149.4
$ws.Range("H137").Value = 48950
$ws.Range("J137").Value = 48950
$ws.Range("L137").Value = 48950
$ws.Range("N137").Value = -59150

